$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Swap LastName / FirstName columns (B and C) for header + all 10 data rows.
for ($r = 1; $r -le 11; $r++) {
    $b = $ws.Cells.Item($r, 2).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 2).Value2 = $c
    $ws.Cells.Item($r, 3).Value2 = $b
}

# 2) Update CNE (column A) numbering: first literal value becomes 21000001,
#    the remaining rows already derive from it via "+1" formulas, so they
#    recalc automatically.
$ws.Range("A2").Value2 = 21000001

# 3) Rename the note headers from Note_AP2x to Note_AP1x (columns D..L).
$noteCols = @(4,5,6,7,8,9,10,11,12)
$n = 1
foreach ($col in $noteCols) {
    $ws.Cells.Item(1, $col).Value2 = "Note_AP1$n"
    $n++
}

# 4) Update the active selection to match the saved state.
$ws.Range("H11").Select()
